$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARN: replace not found: $old"
    }
    return $ok
}

function Get-ParaIndexByText($search) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t -like "*$search*") {
            return $i
        }
    }
    return -1
}

function Insert-BulletAfter($anchorSearch, $newText) {
    $idx = Get-ParaIndexByText($anchorSearch)
    if ($idx -eq -1) {
        Write-Host "WARN: anchor not found: $anchorSearch"
        return
    }
    $anchor = $d.Paragraphs($idx)
    $r = $anchor.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newp = $d.Paragraphs($idx + 1)
    $newp.Range.Text = $newText
}

# ---------------------------------------------------------------------------
# 1. Professional summary
# ---------------------------------------------------------------------------
Replace-Text `
    "Innovative tech leader with over eight years of experience in cloud development, specializing in AI technologies and their application in media contexts. Proven track record of collaborating directly with product and technological teams to develop AI-driven solutions that enhance editorial workflows and consumer-facing products. Skilled in multi-agent orchestration, extensive computing systems, and bridging journalistic principles with AI solutions. Strong advocate for responsible AI usage with a strategic mindset and excellent leadership skills, capable of communicating complex AI concepts to non-technical stakeholders." `
    "Seasoned tech lead with over eight years of experience in cloud development, specializing in AI technologies and their applications in media contexts. Proven track record of collaborating with product and technological teams to develop innovative solutions, with a strong focus on editorial settings. Skilled in multi-agent orchestration, data analysis, and machine learning applications, with a deep understanding of consumer product development. Experienced in senior-level editorial roles, bridging journalistic principles with AI solutions, and advocating for responsible AI usage. Holds a degree in Journalism, Data Science, or equivalent work experience, and excels in leadership, strategic decision-making, and cross-functional collaboration."

# ---------------------------------------------------------------------------
# 2. Role 1 - Global Cloud Inc., Senior Software Development Engineer
# ---------------------------------------------------------------------------
Replace-Text `
    "Led AI-driven infrastructure innovations, enhancing global deployment strategies and operational efficiencies." `
    "Led AI-driven initiatives enhancing global infrastructure, focusing on innovation and cross-team collaboration."

Replace-Text `
    "• Architected batch compute systems, preventing outages and ensuring 99.99% update coverage globally by ?? %." `
    "• Architected batch compute systems for next-gen SDP, helping CrowdStrike avoid another historic outage by ?? %."

Replace-Text `
    "• Spearheaded data warehouse efforts, boosting rollout visibility for leadership across ?? global regions." `
    "• Achieved 99.99% update coverage for hybrid clouds serving RedRock, ClosedAI, Walnut, and BinaryDance by ?? %."

Replace-Text `
    "• Enhanced LLM reasoning with RL self-play, advancing infra-rollout agents for ?? major deployments." `
    "• Spearheaded cross-team efforts building data warehouses, ensuring global rollout visibility for leadership by ?? %."

Insert-BulletAfter `
    "Spearheaded cross-team efforts building data warehouses, ensuring global rollout visibility for leadership by ?? %." `
    "• Pioneered enhancing small LLMs' reasoning via RL self-play and MCTS to ?? build an infra-rollout agent."

# ---------------------------------------------------------------------------
# 3. Role 2 - Global Cloud Inc., Software Development Engineer II
# ---------------------------------------------------------------------------
Replace-Text `
    "Orchestrated large-scale rollouts, optimizing cluster operations and enhancing system reliability." `
    "Developed AI solutions for cluster orchestration, optimizing rollout efficiency and reliability across 8M nodes."

Replace-Text `
    "• Identified 49 metrics for rollout analysis, saving $2M in losses and improving system resilience by ?? %." `
    "• Identified 49 new metrics for rollout failure analysis, saving over $2M in losses for Walnut and RedRock by ?? %."

Replace-Text `
    "• Refined legacy algorithms, cutting rollout time by ?? 30% and improving customer satisfaction." `
    "• Simplified watermark of ?? data aggregation for RTB Ad Exchange, saving advertisers 5% ad revenue loss."

Replace-Text `
    "• Innovated alert merging, reducing Mean Time to ?? Detect from 24h to 30m, ensuring 99.9% SLA uptime." `
    "• Improved customer satisfaction by ?? refining legacy rollout algorithms, cutting rollout time by 30%."

Insert-BulletAfter `
    "Improved customer satisfaction by ?? refining legacy rollout algorithms, cutting rollout time by 30%." `
    "• Innovated alert merging, cutting Mean Time to ?? Detect from 24h to 30m, ensuring 99.9% SLA uptime."

# ---------------------------------------------------------------------------
# 4. Role 3 - TechCorp LLC, Software Engineer
# ---------------------------------------------------------------------------
Replace-Text `
    "Drove ML infrastructure migration, enhancing service availability and cost efficiency." `
    "Engineered scalable AI systems for data-driven platforms, enhancing service availability and cost efficiency."

Replace-Text `
    "• Led ML infra migration to ?? AWS, achieving 99.9% availability for benefit recommendation services." `
    "• Led ML infra migration to ?? AWS, achieving 99.9% availability for benefit recommendation service."

Replace-Text `
    "• Designed distributed message queues, streamlining enterprise integration for ?? major clients." `
    "• Designed a distributed message queue streamlining enterprise integration between BDP and OfficeDay by ?? %."

Replace-Text `
    "• Implemented caching service for breach detection, saving $200k in cloud costs annually by ?? %." `
    "• Implemented a Bloom-Filter caching service for password breach detection, saving over $200k in cloud costs by ?? %."

# ---------------------------------------------------------------------------
# 5. Role 4 - HealthData Systems, Software Engineer
# ---------------------------------------------------------------------------
Replace-Text `
    "Developed cloud solutions for data durability and efficient storage management in healthcare." `
    "Developed cloud-based solutions for data management, ensuring high durability and efficient resource use."

Replace-Text `
    "• Designed data placement service, ensuring 99.9999% data durability for healthcare data lakes by ?? %." `
    "• Designed a data placement service for S3-like storage, ensuring 99.9999% data durability for data lakes by ?? %."

Replace-Text `
    "• Built a garbage collector, reclaiming space and handling ?? TB of deleted and corrupted data." `
    "• Built a garbage collector to ?? reclaim space via compaction, handling deleted, orphaned, and corrupted data."

# ---------------------------------------------------------------------------
# 6. Education - M.S. in Computer Science bullets
# ---------------------------------------------------------------------------
Replace-Text `
    "• Completed advanced coursework in AI technologies and data science, focusing on machine learning applications and data analysis." `
    "• Specialized coursework in AI technologies and data science methodologies"

# Remove the "Led a capstone project..." paragraph entirely.
$idxRemove = Get-ParaIndexByText("Led a capstone project on developing AI-driven solutions for media content analysis, enhancing editorial workflows.")
if ($idxRemove -ne -1) {
    $d.Paragraphs($idxRemove).Range.Delete()
} else {
    Write-Host "WARN: paragraph to delete not found"
}

Replace-Text `
    "• Collaborated with cross-functional teams on projects integrating AI with consumer product development." `
    "• Completed a capstone project on machine learning applications in media, focusing on enhancing editorial workflows"

Insert-BulletAfter `
    "Completed a capstone project on machine learning applications in media, focusing on enhancing editorial workflows" `
    "• Led a research team in developing AI-driven solutions for consumer product development"

# ---------------------------------------------------------------------------
# 7. Education - B.S. in Computer Science bullets
# ---------------------------------------------------------------------------
Replace-Text `
    "• Graduated with honors, specializing in data analysis and machine learning methodologies." `
    "• Graduated with honors, focusing on data analysis and machine learning"

Replace-Text `
    "• Conducted research on AI implications in media, contributing to a publication on AI ethics in journalism." `
    "• Conducted a thesis on the implications of AI in media, exploring editorial and consumer-facing applications"

Replace-Text `
    "• Participated in a collaborative project with the journalism department to develop AI tools for editorial decision-making." `
    "• Participated in a cross-functional collaboration project with the journalism department to bridge AI solutions with journalistic principles"

# ---------------------------------------------------------------------------
# 8. Skills - OTHER list
# ---------------------------------------------------------------------------
Replace-Text `
    "Collaboration with product and technological teams, Editorial operations" `
    "Collaboration with product and technological teams, Editorial operations, Decision-making, Organization"

# ---------------------------------------------------------------------------
# 9. Style spacing: zero-out "space after" for the MR_* paragraph styles
# ---------------------------------------------------------------------------
$styleNames = @("MR_Content", "MR_RoleDescription", "MR_BulletPoint", "MR_SummaryText", "MR_SkillCategory", "MR_SkillList")
foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $s.ParagraphFormat.SpaceAfter = 0
}

Write-Host "All edits applied."
